$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C becomes a static value (was =A+B formula)
$ws.Range("C1").Value = 3
$ws.Range("C2").Value = 3
$ws.Range("C3").Value = 5
$ws.Range("C4").Value = 4

# Insert new column D (shifts old D -> E) and fill with A*B*C formula
$ws.Range("D1:D4").Insert(-4161)
$ws.Range("D1").Formula = "=A1*B1*C1"
$ws.Range("D2").Formula = "=A2*B2*C2"
$ws.Range("D3").Formula = "=A3*B3*C3"
$ws.Range("D4").Formula = "=A4*B4*C4"

# Columns F and G: lookup-style text values
$ws.Range("F1").Value = "sample1"
$ws.Range("G1").Value = "sample2"

$ws.Range("F2").Value = "sample2"
$ws.Range("G2").Value = "sample2"

$ws.Range("F3").Value = "sample3"
$ws.Range("G3").Value = "sample3"

$ws.Range("F4").Value = "sample1"
$ws.Range("G4").Value = "sample1"

# Column H: SUMPRODUCT formula
$ws.Range("H1").Formula = "=SUMPRODUCT((F`$1:F`$9519=F1)*(G`$1:G`$9519=G1)*C`$2:C`$9519)"
$ws.Range("H2").Formula = "=SUMPRODUCT((F`$1:F`$9519=F2)*(G`$1:G`$9519=G2)*C`$2:C`$9519)"
$ws.Range("H3").Formula = "=SUMPRODUCT((F`$1:F`$9519=F3)*(G`$1:G`$9519=G3)*C`$2:C`$9519)"
$ws.Range("H4").Formula = "=SUMPRODUCT((F`$1:F`$9519=F4)*(G`$1:G`$9519=G4)*C`$2:C`$9519)"

$ws.Range("H1").Select()
